# Update last 30 days report (1/12/2025)
# - Insert a new tracking row for the Dash0 / Director of Commercial Sales role
#   (Tim Sullivan moving to 2nd Interview), shifting all subsequent rows down by one.
# - Update the existing Dash0 row (row 7, Chris Hogan) with a refreshed action date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8; this shifts current rows 8-23 down to 9-24
# and duplicates formatting (incl. the date style on column F) from the row above.
$ws.Rows.Item(8).Insert()

# Update row 7 (Dash0 / Chris Hogan) - candidate unchanged company/job, date + note refreshed
$ws.Cells.Item(7, 4).Value = "Chris Hogan"
$ws.Cells.Item(7, 6).Value = 45994

# Populate the newly inserted row 8 (Dash0 / Tim Sullivan moves to 2nd Interview)
$ws.Cells.Item(8, 1).Value = 738
$ws.Cells.Item(8, 2).Value = "Dash0"
$ws.Cells.Item(8, 3).Value = "Director of Commercial Sales, NYC/Boston"
$ws.Cells.Item(8, 4).Value = "Tim Sullivan"
$ws.Cells.Item(8, 5).Value = "2nd Interview"
$ws.Cells.Item(8, 6).Value = 45994
